$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (source data is inline strings;
# some values like "206.64" would otherwise auto-convert to numbers).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.919.60'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.552.22'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '206.64'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.487'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.51%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +1.73%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.248'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0857'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.772.77'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.544.02'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.06%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '26.908.53'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '61.67'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '217.12'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0689'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.24'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.03%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.86'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '6.58'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '14.86'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.421.57'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.09%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +4.20%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.957'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +1.98%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.71%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.521'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.985'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '63.68'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.687.64'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '86.18'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.82%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0₇0991'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.79%  '
